$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("remittances")

# Column D: "category group" labels for the remittances sheet (new draft for voronoi treemap)
$groups = @(
    "Food",
    "Food",
    "Utilities",
    "Housing",
    "Utilities",
    "Housing",
    "Transportation",
    "Utilities",
    "Utilities",
    "Personal Spending",
    "Healthcare",
    "Personal Spending",
    "Education",
    "Savings & Debt",
    "Personal Spending",
    "Productive Supplies",
    "Savings & Debt",
    "Housing"
)

for ($i = 0; $i -lt $groups.Length; $i++) {
    $row = $i + 1
    $ws1.Range("D$row").Value = $groups[$i]
}

# Make "remittances" the active sheet/tab, with C10 selected,
# which also clears the previous tabSelected/active state on "No remittances".
$ws1.Activate()
$ws1.Range("C10").Select()
